# Apply the "Fixed update to excel issue" change:
#  1. Rename the "Requested quantity" header on the "Weekly Quantity" sheet to "Weekly_PO_Qty"
#  2. Rename the "Requested quantity" header on the "Monthly Trend" sheet to "Monthly_PO_Qty"
#  3. Add a new "PO Forecast" worksheet (after the existing sheets) with forecast data

$wb = $excel.ActiveWorkbook

# --- 1. Weekly Quantity sheet: rename header ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

# --- 2. Monthly Trend sheet: rename header ---
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- 3. Add new "PO Forecast" sheet at the end ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add($null, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Copy header formatting (bold/centered/bordered) from the Weekly Quantity header row
$wsWeekly.Range("A1:B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)

# Copy the date-column formatting from the Weekly Quantity data column
$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A19").PasteSpecial(-4122)

# Header values
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Forecast data rows
$wsForecast.Range("A2").Value = 45466.99999999999
$wsForecast.Range("B2").Value = 24
$wsForecast.Range("C2").Value = -7.348276675640614
$wsForecast.Range("D2").Value = 54.10497444634085

$wsForecast.Range("A3").Value = 45480.99999999999
$wsForecast.Range("B3").Value = 23
$wsForecast.Range("C3").Value = -5.334052503704056
$wsForecast.Range("D3").Value = 53.50575465154403

$wsForecast.Range("A4").Value = 45487.99999999999
$wsForecast.Range("B4").Value = 23
$wsForecast.Range("C4").Value = -5.533173635785066
$wsForecast.Range("D4").Value = 51.90506325052872

$wsForecast.Range("A5").Value = 45494.99999999999
$wsForecast.Range("B5").Value = 23
$wsForecast.Range("C5").Value = -5.499803312561552
$wsForecast.Range("D5").Value = 51.59088599254965

$wsForecast.Range("A6").Value = 45501.99999999999
$wsForecast.Range("B6").Value = 22
$wsForecast.Range("C6").Value = -8.119018698198648
$wsForecast.Range("D6").Value = 50.76541744928284

$wsForecast.Range("A7").Value = 45508.99999999999
$wsForecast.Range("B7").Value = 22
$wsForecast.Range("C7").Value = -7.441823402412783
$wsForecast.Range("D7").Value = 52.57108295679698

$wsForecast.Range("A8").Value = 45515.99999999999
$wsForecast.Range("B8").Value = 22
$wsForecast.Range("C8").Value = -8.506519760359907
$wsForecast.Range("D8").Value = 48.99621474322637

$wsForecast.Range("A9").Value = 45522.99999999999
$wsForecast.Range("B9").Value = 21
$wsForecast.Range("C9").Value = -8.121877963404321
$wsForecast.Range("D9").Value = 50.44758206203306

$wsForecast.Range("A10").Value = 45634.99999999999
$wsForecast.Range("B10").Value = 16
$wsForecast.Range("C10").Value = -14.02335514208442
$wsForecast.Range("D10").Value = 46.37107757001373

$wsForecast.Range("A11").Value = 45641.99999999999
$wsForecast.Range("B11").Value = 16
$wsForecast.Range("C11").Value = -12.72512607497063
$wsForecast.Range("D11").Value = 46.29706292727251

$wsForecast.Range("A12").Value = 45648.99999999999
$wsForecast.Range("B12").Value = 16
$wsForecast.Range("C12").Value = -14.14037799978814
$wsForecast.Range("D12").Value = 43.98378513301252

$wsForecast.Range("A13").Value = 45655.99999999999
$wsForecast.Range("B13").Value = 16
$wsForecast.Range("C13").Value = -13.69550124038839
$wsForecast.Range("D13").Value = 46.85071396642169

$wsForecast.Range("A14").Value = 45662.99999999999
$wsForecast.Range("B14").Value = 15
$wsForecast.Range("C14").Value = -14.68215864486471
$wsForecast.Range("D14").Value = 44.14766909507512

$wsForecast.Range("A15").Value = 45669.99999999999
$wsForecast.Range("B15").Value = 15
$wsForecast.Range("C15").Value = -14.53830136740547
$wsForecast.Range("D15").Value = 43.61046220013424

$wsForecast.Range("A16").Value = 45676.99999999999
$wsForecast.Range("B16").Value = 15
$wsForecast.Range("C16").Value = -16.60132562284379
$wsForecast.Range("D16").Value = 43.66795149101696

$wsForecast.Range("A17").Value = 45683.99999999999
$wsForecast.Range("B17").Value = 14
$wsForecast.Range("C17").Value = -15.4143203917323
$wsForecast.Range("D17").Value = 45.79724100858104

$wsForecast.Range("A18").Value = 45690.99999999999
$wsForecast.Range("B18").Value = 14
$wsForecast.Range("C18").Value = -15.6778996858061
$wsForecast.Range("D18").Value = 44.4839552913019

$wsForecast.Range("A19").Value = 45697.99999999999
$wsForecast.Range("B19").Value = 14
$wsForecast.Range("C19").Value = -15.20896591693077
$wsForecast.Range("D19").Value = 43.20780405583471

$wsForecast.Range("A1").Select()
